$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2Player1")

# First, copy row 7's current (custom AM/PM) duration format down onto row
# 10 -- that's the row that ends up carrying it in the final layout -- while
# the source cell still has it, before row 7 gets reused for other data.
$ws.Range("C7").Copy($ws.Range("C10"))

# NOTE: new shared strings must be interned in this exact order so they
# land at shared-string table indices 30,31,32,33 respectively, matching
# the target workbook (test -> dsfsdgfs -> sdfgsdf -> asdf).
$ws.Cells.Item(6, 1).Value = "test"
$ws.Cells.Item(13, 1).Value = "dsfsdgfs"
$ws.Cells.Item(4, 1).Value = "sdfgsdf"
$ws.Cells.Item(5, 1).Value = "asdf"

$data = @(
    @{Row=3;  Name="Tjeerd";   Score=400; Seconds=2},
    @{Row=4;  Name="sdfgsdf";  Score=400; Seconds=2},
    @{Row=5;  Name="asdf";     Score=400; Seconds=4},
    @{Row=6;  Name="test";     Score=300; Seconds=10},
    @{Row=7;  Name="Testasda"; Score=250; Seconds=1},
    @{Row=8;  Name="Tjeerd";   Score=200; Seconds=1},
    @{Row=9;  Name="Tjeerd";   Score=200; Seconds=1},
    @{Row=10; Name="Test";     Score=150; Seconds=1},
    @{Row=11; Name="test";     Score=100; Seconds=10},
    @{Row=12; Name="1sdfsdfs"; Score=0;   Seconds=2},
    @{Row=13; Name="dsfsdgfs"; Score=0;   Seconds=2},
    @{Row=14; Name="asdf";     Score=0;   Seconds=4}
)

foreach ($row in $data) {
    $ws.Cells.Item($row.Row, 1).Value = $row.Name
    $ws.Cells.Item($row.Row, 2).Value = $row.Score
    $ws.Cells.Item($row.Row, 3).Value = (1 / 86400) * $row.Seconds
}

# Row 7 now holds an ordinary duration entry, same plain format as the
# rest of the table (its old custom AM/PM format moved to row 10 above).
$ws.Range("C3").Copy($ws.Range("C7"))
$ws.Cells.Item(7, 3).Value = (1 / 86400) * 1

# Rows 11-14 are brand new table rows (beyond the sheet's previous used
# range) -- give their duration cells the same plain format as the rest
# of the column.
$ws.Range("C3").Copy($ws.Range("C11:C14"))
$ws.Cells.Item(11, 3).Value = (1 / 86400) * 10
$ws.Cells.Item(12, 3).Value = (1 / 86400) * 2
$ws.Cells.Item(13, 3).Value = (1 / 86400) * 2
$ws.Cells.Item(14, 3).Value = (1 / 86400) * 4

# Refresh the worksheet's remembered "last sort" state (range + key +
# direction) to cover the now-larger table, matching how Excel records
# sortState after a Data > Sort on A3:C14 keyed on column B descending.
$sort = $ws.Sort
$sortFields = $sort.SortFields
$sortFields.Clear()
$sortFields.Add($ws.Range("B3"), 0, 2, $null, 0)
$sort.SetRange($ws.Range("A3:C14"))
$sort.Header = 2
$sort.Apply()
